$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "Gestionar las calificaciones"
